$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Vegfc"
$ws.Range("C2").Value = "Kdr"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.663313
$ws.Range("H2").Value = 7.989939000000001
$ws.Range("I2").Value = 0.3794306644527501
$ws.Range("J2").Value = 0.3794306644527502
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 139.3946303333333
$ws.Range("N2").Value = 418.183891
$ws.Range("O2").Value = 0.9207771771472824
$ws.Range("P2").Value = 0.9207771771472822
$ws.Range("Q2").Value = 371.2515310969609
$ws.Range("R2").Value = 3341.263779872649
$ws.Range("S2").Value = 0.3493710961379209
$ws.Range("T2").Value = 0.349371096137921

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Vegfc"
$ws.Range("C3").Value = "Kdr"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.663313
$ws.Range("H3").Value = 7.989939000000001
$ws.Range("I3").Value = 0.3794306644527501
$ws.Range("J3").Value = 0.3794306644527502
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.050239
$ws.Range("N3").Value = 0.150717
$ws.Range("O3").Value = 0.0003318558576616883
$ws.Range("P3").Value = 0.0003318558576616882
$ws.Range("Q3").Value = 0.133802181807
$ws.Range("R3").Value = 1.204219636263
$ws.Range("S3").Value = 0.0001259162885751117
$ws.Range("T3").Value = 0.0001259162885751117

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Vegfc"
$ws.Range("C4").Value = "Kdr"
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.663313
$ws.Range("H4").Value = 7.989939000000001
$ws.Range("I4").Value = 0.3794306644527501
$ws.Range("J4").Value = 0.3794306644527502
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 4.220699666666667
$ws.Range("N4").Value = 12.662099
$ws.Range("O4").Value = 0.02788001170035368
$ws.Range("P4").Value = 0.02788001170035368
$ws.Range("Q4").Value = 11.241044291329
$ws.Range("R4").Value = 101.169398621961
$ws.Range("S4").Value = 0.01057853136441565
$ws.Range("T4").Value = 0.01057853136441565

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Vegfc"
$ws.Range("C5").Value = "Kdr"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.663313
$ws.Range("H5").Value = 7.989939000000001
$ws.Range("I5").Value = 0.3794306644527501
$ws.Range("J5").Value = 0.3794306644527502
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 6.099343
$ws.Range("N5").Value = 18.298029
$ws.Range("O5").Value = 0.04028947038033828
$ws.Range("P5").Value = 0.04028947038033828
$ws.Range("Q5").Value = 16.244459503359
$ws.Range("R5").Value = 146.200135530231
$ws.Range("S5").Value = 0.01528706051686115
$ws.Range("T5").Value = 0.01528706051686115

# Row 6
$ws.Range("A6").Value = "ECs"
$ws.Range("B6").Value = "Vegfc"
$ws.Range("C6").Value = "Kdr"
$ws.Range("D6").Value = "Neutro"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 2.663313
$ws.Range("H6").Value = 7.989939000000001
$ws.Range("I6").Value = 0.3794306644527501
$ws.Range("J6").Value = 0.3794306644527502
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.8278236666666666
$ws.Range("N6").Value = 2.483471
$ws.Range("O6").Value = 0.005468224544563193
$ws.Range("P6").Value = 0.005468224544563191
$ws.Range("Q6").Value = 2.204753533141
$ws.Range("R6").Value = 19.842781798269
$ws.Range("S6").Value = 0.002074812072320449
$ws.Range("T6").Value = 0.002074812072320449

# Row 7
$ws.Range("A7").Value = "ECs"
$ws.Range("B7").Value = "Vegfc"
$ws.Range("C7").Value = "Kdr"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 2.663313
$ws.Range("H7").Value = 7.989939000000001
$ws.Range("I7").Value = 0.3794306644527501
$ws.Range("J7").Value = 0.3794306644527502
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.7952806666666667
$ws.Range("N7").Value = 2.385842
$ws.Range("O7").Value = 0.005253260369800871
$ws.Range("P7").Value = 0.00525326036980087
$ws.Range("Q7").Value = 2.118081338182
$ws.Range("R7").Value = 19.062732043638
$ws.Range("S7").Value = 0.001993248072656844
$ws.Range("T7").Value = 0.001993248072656844

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Vegfc"
$ws.Range("C8").Value = "Kdr"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 3.178631333333334
$ws.Range("H8").Value = 9.535894000000001
$ws.Range("I8").Value = 0.4528458348143826
$ws.Range("J8").Value = 0.4528458348143827
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 139.3946303333333
$ws.Range("N8").Value = 418.183891
$ws.Range("O8").Value = 0.9207771771472824
$ws.Range("P8").Value = 0.9207771771472822
$ws.Range("Q8").Value = 443.0841396759504
$ws.Range("R8").Value = 3987.757257083554
$ws.Range("S8").Value = 0.4169701094632918
$ws.Range("T8").Value = 0.4169701094632918

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Vegfc"
$ws.Range("C9").Value = "Kdr"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 3.178631333333334
$ws.Range("H9").Value = 9.535894000000001
$ws.Range("I9").Value = 0.4528458348143826
$ws.Range("J9").Value = 0.4528458348143827
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.050239
$ws.Range("N9").Value = 0.150717
$ws.Range("O9").Value = 0.0003318558576616883
$ws.Range("P9").Value = 0.0003318558576616882
$ws.Range("Q9").Value = 0.1596912595553333
$ws.Range("R9").Value = 1.437221335998
$ws.Range("S9").Value = 0.0001502795429008502
$ws.Range("T9").Value = 0.0001502795429008502

# Row 10
$ws.Range("A10").Value = "FAPs"
$ws.Range("B10").Value = "Vegfc"
$ws.Range("C10").Value = "Kdr"
$ws.Range("D10").Value = "M1"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 3.178631333333334
$ws.Range("H10").Value = 9.535894000000001
$ws.Range("I10").Value = 0.4528458348143826
$ws.Range("J10").Value = 0.4528458348143827
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 4.220699666666667
$ws.Range("N10").Value = 12.662099
$ws.Range("O10").Value = 0.02788001170035368
$ws.Range("P10").Value = 0.02788001170035368
$ws.Range("Q10").Value = 13.41604820905622
$ws.Range("R10").Value = 120.744433881506
$ws.Range("S10").Value = 0.01262534717308142
$ws.Range("T10").Value = 0.01262534717308142

# Row 11
$ws.Range("A11").Value = "FAPs"
$ws.Range("B11").Value = "Vegfc"
$ws.Range("C11").Value = "Kdr"
$ws.Range("D11").Value = "M2"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 3.178631333333334
$ws.Range("H11").Value = 9.535894000000001
$ws.Range("I11").Value = 0.4528458348143826
$ws.Range("J11").Value = 0.4528458348143827
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 6.099343
$ws.Range("N11").Value = 18.298029
$ws.Range("O11").Value = 0.04028947038033828
$ws.Range("P11").Value = 0.04028947038033828
$ws.Range("Q11").Value = 19.38756277254734
$ws.Range("R11").Value = 174.488064952926
$ws.Range("S11").Value = 0.01824491884861363
$ws.Range("T11").Value = 0.01824491884861363

# Row 12
$ws.Range("A12").Value = "FAPs"
$ws.Range("B12").Value = "Vegfc"
$ws.Range("C12").Value = "Kdr"
$ws.Range("D12").Value = "Neutro"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 3.178631333333334
$ws.Range("H12").Value = 9.535894000000001
$ws.Range("I12").Value = 0.4528458348143826
$ws.Range("J12").Value = 0.4528458348143827
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.8278236666666666
$ws.Range("N12").Value = 2.483471
$ws.Range("O12").Value = 0.005468224544563193
$ws.Range("P12").Value = 0.005468224544563191
$ws.Range("Q12").Value = 2.631346245341556
$ws.Range("R12").Value = 23.682116208074
$ws.Range("S12").Value = 0.002476262708835216
$ws.Range("T12").Value = 0.002476262708835216

# Row 13
$ws.Range("A13").Value = "FAPs"
$ws.Range("B13").Value = "Vegfc"
$ws.Range("C13").Value = "Kdr"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 3.178631333333334
$ws.Range("H13").Value = 9.535894000000001
$ws.Range("I13").Value = 0.4528458348143826
$ws.Range("J13").Value = 0.4528458348143827
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.7952806666666667
$ws.Range("N13").Value = 2.385842
$ws.Range("O13").Value = 0.005253260369800871
$ws.Range("P13").Value = 0.00525326036980087
$ws.Range("Q13").Value = 2.527904045860889
$ws.Range("R13").Value = 22.751136412748
$ws.Range("S13").Value = 0.002378917077659787
$ws.Range("T13").Value = 0.002378917077659787

# Row 14
$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Vegfc"
$ws.Range("C14").Value = "Kdr"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 1.177290666666667
$ws.Range("H14").Value = 3.531872
$ws.Range("I14").Value = 0.1677235007328671
$ws.Range("J14").Value = 0.1677235007328671
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 139.3946303333333
$ws.Range("N14").Value = 418.183891
$ws.Range("O14").Value = 0.9207771771472824
$ws.Range("P14").Value = 0.9207771771472822
$ws.Range("Q14").Value = 164.1079972748835
$ws.Range("R14").Value = 1476.971975473952
$ws.Range("S14").Value = 0.1544359715460695
$ws.Range("T14").Value = 0.1544359715460695

# Row 15
$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Vegfc"
$ws.Range("C15").Value = "Kdr"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 1.177290666666667
$ws.Range("H15").Value = 3.531872
$ws.Range("I15").Value = 0.1677235007328671
$ws.Range("J15").Value = 0.1677235007328671
$ws.Range("K15").Value = 1
$ws.Range("L15").Value = 0.3333333333333333
$ws.Range("M15").Value = 0.050239
$ws.Range("N15").Value = 0.150717
$ws.Range("O15").Value = 0.0003318558576616883
$ws.Range("P15").Value = 0.0003318558576616882
$ws.Range("Q15").Value = 0.05914590580266667
$ws.Range("R15").Value = 0.532313152224
$ws.Range("S15").Value = 0.00005566002618572642
$ws.Range("T15").Value = 0.00005566002618572642

# Row 16
$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Vegfc"
$ws.Range("C16").Value = "Kdr"
$ws.Range("D16").Value = "M1"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 1.177290666666667
$ws.Range("H16").Value = 3.531872
$ws.Range("I16").Value = 0.1677235007328671
$ws.Range("J16").Value = 0.1677235007328671
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 4.220699666666667
$ws.Range("N16").Value = 12.662099
$ws.Range("O16").Value = 0.02788001170035368
$ws.Range("P16").Value = 0.02788001170035368
$ws.Range("Q16").Value = 4.968990324369778
$ws.Range("R16").Value = 44.720912919328
$ws.Range("S16").Value = 0.004676133162856614
$ws.Range("T16").Value = 0.004676133162856614

# Row 17
$ws.Range("A17").Value = "sCs"
$ws.Range("B17").Value = "Vegfc"
$ws.Range("C17").Value = "Kdr"
$ws.Range("D17").Value = "M2"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 1.177290666666667
$ws.Range("H17").Value = 3.531872
$ws.Range("I17").Value = 0.1677235007328671
$ws.Range("J17").Value = 0.1677235007328671
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 6.099343
$ws.Range("N17").Value = 18.298029
$ws.Range("O17").Value = 0.04028947038033828
$ws.Range("P17").Value = 0.04028947038033828
$ws.Range("Q17").Value = 7.180699586698667
$ws.Range("R17").Value = 64.626296280288
$ws.Range("S17").Value = 0.006757491014863496
$ws.Range("T17").Value = 0.006757491014863495

# Row 18
$ws.Range("A18").Value = "sCs"
$ws.Range("B18").Value = "Vegfc"
$ws.Range("C18").Value = "Kdr"
$ws.Range("D18").Value = "Neutro"
$ws.Range("E18").Value = 3
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 1.177290666666667
$ws.Range("H18").Value = 3.531872
$ws.Range("I18").Value = 0.1677235007328671
$ws.Range("J18").Value = 0.1677235007328671
$ws.Range("K18").Value = 3
$ws.Range("L18").Value = 1
$ws.Range("M18").Value = 0.8278236666666666
$ws.Range("N18").Value = 2.483471
$ws.Range("O18").Value = 0.005468224544563193
$ws.Range("P18").Value = 0.005468224544563191
$ws.Range("Q18").Value = 0.9745890764124444
$ws.Range("R18").Value = 8.771301687712
$ws.Range("S18").Value = 0.0009171497634075267
$ws.Range("T18").Value = 0.0009171497634075265

# Row 19
$ws.Range("A19").Value = "sCs"
$ws.Range("B19").Value = "Vegfc"
$ws.Range("C19").Value = "Kdr"
$ws.Range("D19").Value = "sCs"
$ws.Range("E19").Value = 3
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 1.177290666666667
$ws.Range("H19").Value = 3.531872
$ws.Range("I19").Value = 0.1677235007328671
$ws.Range("J19").Value = 0.1677235007328671
$ws.Range("K19").Value = 3
$ws.Range("L19").Value = 1
$ws.Range("M19").Value = 0.7952806666666667
$ws.Range("N19").Value = 2.385842
$ws.Range("O19").Value = 0.005253260369800871
$ws.Range("P19").Value = 0.00525326036980087
$ws.Range("Q19").Value = 0.9362765062471112
$ws.Range("R19").Value = 8.426488556224001
$ws.Range("S19").Value = 0.0008810952194842381
$ws.Range("T19").Value = 0.000881095219484238
